# Updates the EC (Estado de Cuenta) workbook:
#  - Removes the 5 rows of old worker data (MARGOTH LUCIA HERNANDEZ ALVAREZ,
#    periods 2507-2503), shifting the remaining worker rows (ZANIA, JUAN CARLOS)
#    and the signature block up.
#  - Refreshes the "Valor Mora", "Cant. Trabajadores" and "Cant. Periodos" totals.
#  - Narrows column D slightly to fit the new (shorter) longest name.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the 5 rows belonging to MARGOTH LUCIA HERNANDEZ ALVAREZ (rows 16-20).
# This shifts the remaining data rows (ZANIA, JUAN CARLOS) and the signature
# rows below them up by 5 rows automatically, along with their merged cells.
$ws.Range("A16:A20").EntireRow.Delete()

# Update the summary figures at the top of the statement.
$ws.Range("E11").Value = 61967
$ws.Range("C13").Value = 2
$ws.Range("F13").Value = 2

# Column D is now narrower since the remaining names are shorter.
# (ColumnWidth uses the default-font character unit; 33.15 round-trips to a
# stored column width of exactly 34 in the saved OOXML.)
$ws.Columns.Item(4).ColumnWidth = 33.15
